$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '60.247.90'
$ws.Range('E2').Value = '  +1.68%  '
$ws.Range('D3').Value = '2.680.78'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '522.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('E6').Value = '  -0.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('E8').Value = '  +0.88%  '
$ws.Range('D9').Value = '2.700.52'
$ws.Range('E9').Value = '  -1.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.47'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.03%  '
$ws.Range('E11').Value = '  -1.73%  '
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('E13').Value = '  +1.56%  '
$ws.Range('D14').Value = '3.156.62'
$ws.Range('E14').Value = '  -0.26%  '
$ws.Range('D15').Value = '60.252.91'
$ws.Range('E15').Value = '  +1.79%  '
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('D18').Value = '2.683.23'
$ws.Range('E18').Value = '  -1.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '351.31'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('E20').Value = '  -0.89%  '
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('E22').Value = '  +1.38%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.15'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.90%  '
$ws.Range('E25').Value = '  -0.86%  '
$ws.Range('E26').Value = '  +4.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.993'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.25%  '
$ws.Range('E28').Value = '  +0.83%  '
$ws.Range('D29').Value = '0.0₃0818'
$ws.Range('E29').Value = '  -1.47%  '
$ws.Range('E30').Value = '  +5.36%  '
$ws.Range('E31').Value = '  +0.30%  '
$ws.Range('E32').Value = '  +0.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.14'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '147.22'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.32'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.37%  '
$ws.Range('E36').Value = '  +8.24%  '
$ws.Range('E37').Value = '  -7.35%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.878'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.39%  '
$ws.Range('E39').Value = '  +6.90%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.94'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.72'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.69%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '284.88'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.07'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0992'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.611'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.80%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.996'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.53%  '
$ws.Range('D47').Value = '2.131.29'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.96'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.41%  '
$ws.Range('E49').Value = '  +0.79%  '
$ws.Range('E50').Value = '  +1.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.40'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.64%  '
